$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Journal entry for "Friday 22 April 2022" (row 17): the end time of the
# session moved from 03:00 to 06:00, so the logged duration grows from
# 10 hours to 13 hours. The activity description is unchanged.
$ws.Range("C17").Value = "17:00-06:00"
$ws.Range("D17").Value = 13

# Move the active selection to D18 (was E9).
$ws.Activate()
$ws.Range("D18").Select()
